$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 160
$ws.Range("I39").Value = 95
$ws.Range("J39").Value = 290
$ws.Range("K39").Value = 285
$ws.Range("L39").Value = 870
$ws.Range("M39").Value = 11
$ws.Range("N39").Value = -1462
$ws.Range("H69").Value = 3351.3635
$ws.Range("I69").Value = 3336.5
$ws.Range("J69").Value = 3500
$ws.Range("K69").Value = 10009.5
$ws.Range("L69").Value = 10500
$ws.Range("M69").Value = -9135.5
$ws.Range("N69").Value = -12248
$ws.Range("H72").Value = 3351.3635
$ws.Range("I72").Value = 3336.5
$ws.Range("J72").Value = 3500
$ws.Range("K72").Value = 30028.5
$ws.Range("L72").Value = 31500
$ws.Range("M72").Value = -25660.5
$ws.Range("N72").Value = -40236
$ws.Range("H100").Value = 1970.3334
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 1970.3334
$ws.Range("K100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("M100").Value = 1970.3334
$ws.Range("N100").Value = -3052.3334
$ws.Range("H107").Value = 13374.9
$ws.Range("I107").Value = 18749.572
$ws.Range("J107").Value = 834
$ws.Range("K107").Value = 18749.572
$ws.Range("L107").Value = 834
$ws.Range("M107").Value = -16829.572
$ws.Range("N107").Value = -4674
$ws.Range("H116").Value = 1929.875
$ws.Range("I116").Value = 1898.3077
$ws.Range("J116").Value = 2066.6667
$ws.Range("K116").Value = 1898.3077
$ws.Range("L116").Value = 2066.6667
$ws.Range("M116").Value = 1543.6923
$ws.Range("N116").Value = -8950.6667
$ws.Range("H123").Value = 26222.5
$ws.Range("J123").Value = 26222.5
$ws.Range("L123").Value = 26222.5
$ws.Range("N123").Value = -36022.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 11556.8
$ws.Range("I34").Value = 4850
$ws.Range("J34").Value = 16028
$ws.Range("K34").Value = 4850
$ws.Range("L34").Value = 16028
$ws.Range("M34").Value = -4579
$ws.Range("N34").Value = -16570
$ws.Range("H74").Value = 15626196
$ws.Range("I74").Value = 19231920
$ws.Range("J74").Value = 1394.6666
$ws.Range("K74").Value = 19231920
$ws.Range("L74").Value = 1394.6666
$ws.Range("M74").Value = -19231046
$ws.Range("N74").Value = -3142.6666
$ws.Range("H77").Value = 15626196
$ws.Range("I77").Value = 19231920
$ws.Range("J77").Value = 1394.6666
$ws.Range("K77").Value = 96159600
$ws.Range("L77").Value = 6973.333000000001
$ws.Range("M77").Value = -96155232
$ws.Range("N77").Value = -15709.333
$ws.Range("H102").Value = 2541.7778
$ws.Range("I102").Value = 2346
$ws.Range("J102").Value = 2933.3333
$ws.Range("K102").Value = 2346
$ws.Range("L102").Value = 2933.3333
$ws.Range("M102").Value = -724
$ws.Range("N102").Value = -6177.3333
$ws.Range("H122").Value = 9638.92
$ws.Range("I122").Value = 9172.739
$ws.Range("K122").Value = 27518.217
$ws.Range("M122").Value = -25068.217

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 39263.332
$ws.Range("J35").Value = 39263.332
$ws.Range("L35").Value = 39263.332
$ws.Range("N35").Value = -39883.332
$ws.Range("H105").Value = 2919.6345
$ws.Range("I105").Value = 1569.7587
$ws.Range("J105").Value = 4621.6523
$ws.Range("K105").Value = 1569.7587
$ws.Range("L105").Value = 4621.6523
$ws.Range("M105").Value = 177.2412999999999
$ws.Range("N105").Value = -8115.6523

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2532.5
$ws.Range("I62").Value = 2214.2856
$ws.Range("K62").Value = 2214.2856
$ws.Range("M62").Value = -1590.2856
$ws.Range("H65").Value = 2532.5
$ws.Range("I65").Value = 2214.2856
$ws.Range("K65").Value = 11071.428
$ws.Range("M65").Value = -7951.428
$ws.Range("H86").Value = 4155.778
$ws.Range("I86").Value = 4902.3335
$ws.Range("J86").Value = 3782.5
$ws.Range("K86").Value = 4902.3335
$ws.Range("L86").Value = 3782.5
$ws.Range("M86").Value = -3779.3335
$ws.Range("N86").Value = -6028.5
$ws.Range("H89").Value = 4155.778
$ws.Range("I89").Value = 4902.3335
$ws.Range("J89").Value = 3782.5
$ws.Range("K89").Value = 24511.6675
$ws.Range("L89").Value = 18912.5
$ws.Range("M89").Value = -18895.6675
$ws.Range("N89").Value = -30144.5
$ws.Range("H105").Value = 1247.8948
$ws.Range("I105").Value = 1016.1539
$ws.Range("J105").Value = 1750
$ws.Range("K105").Value = 1016.1539
$ws.Range("L105").Value = 1750
$ws.Range("M105").Value = 730.8461
$ws.Range("N105").Value = -5244
$ws.Range("H109").Value = 25267.777
$ws.Range("J109").Value = 25926.25
$ws.Range("L109").Value = 25926.25
$ws.Range("N109").Value = -28006.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1503.1781
$ws.Range("I131").Value = 2290.5
$ws.Range("J131").Value = 1206.0754
$ws.Range("K131").Value = 6871.5
$ws.Range("L131").Value = 3618.2262
$ws.Range("M131").Value = -1831.5
$ws.Range("N131").Value = -13698.2262

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 14484.875
$ws.Range("J57").Value = 20996
$ws.Range("L57").Value = 20996
$ws.Range("N57").Value = -22636
$ws.Range("H80").Value = 18521430
$ws.Range("I80").Value = 41669190
$ws.Range("J80").Value = 3218
$ws.Range("K80").Value = 41669190
$ws.Range("L80").Value = 3218
$ws.Range("M80").Value = -41668192
$ws.Range("N80").Value = -5214
$ws.Range("H83").Value = 18521430
$ws.Range("I83").Value = 41669190
$ws.Range("J83").Value = 3218
$ws.Range("K83").Value = 208345950
$ws.Range("L83").Value = 16090
$ws.Range("M83").Value = -208340958
$ws.Range("N83").Value = -26074
$ws.Range("H113").Value = 63647.75
$ws.Range("I113").Value = 125687.625
$ws.Range("J113").Value = 1607.875
$ws.Range("K113").Value = 125687.625
$ws.Range("L113").Value = 1607.875
$ws.Range("M113").Value = -123517.625
$ws.Range("N113").Value = -5947.875
$ws.Range("H132").Value = 6294.5293
$ws.Range("I132").Value = 7010.76
$ws.Range("J132").Value = 4305
$ws.Range("K132").Value = 21032.28
$ws.Range("L132").Value = 12915
$ws.Range("M132").Value = -18502.28
$ws.Range("N132").Value = -17975

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6307.2
$ws.Range("I7").Value = 10734.5
$ws.Range("J7").Value = 4909.1055
$ws.Range("K7").Value = 10734.5
$ws.Range("L7").Value = 4909.1055
$ws.Range("M7").Value = -10622.5
$ws.Range("N7").Value = -5133.1055
$ws.Range("H40").Value = 3641.9333
$ws.Range("I40").Value = 5205.3335
$ws.Range("J40").Value = 2599.6667
$ws.Range("K40").Value = 5205.3335
$ws.Range("L40").Value = 2599.6667
$ws.Range("M40").Value = -5069.3335
$ws.Range("N40").Value = -2871.6667
$ws.Range("H68").Value = 2151.9285
$ws.Range("I68").Value = 2245.5
$ws.Range("J68").Value = 2027.1666
$ws.Range("K68").Value = 2245.5
$ws.Range("L68").Value = 2027.1666
$ws.Range("M68").Value = -1496.5
$ws.Range("N68").Value = -3525.1666
$ws.Range("H71").Value = 2151.9285
$ws.Range("I71").Value = 2245.5
$ws.Range("J71").Value = 2027.1666
$ws.Range("K71").Value = 11227.5
$ws.Range("L71").Value = 10135.833
$ws.Range("M71").Value = -7483.5
$ws.Range("N71").Value = -17623.833
$ws.Range("H94").Value = 46580
$ws.Range("J94").Value = 46580
$ws.Range("L94").Value = 46580
$ws.Range("N94").Value = -47932
$ws.Range("H100").Value = 2330.4
$ws.Range("I100").Value = 1956
$ws.Range("K100").Value = 1956
$ws.Range("M100").Value = -1415
$ws.Range("H122").Value = 6285.879
$ws.Range("I122").Value = 6824.0557
$ws.Range("J122").Value = 5640.067
$ws.Range("K122").Value = 20472.1671
$ws.Range("L122").Value = 16920.201
$ws.Range("M122").Value = -18022.1671
$ws.Range("N122").Value = -21820.201
$ws.Range("H126").Value = 6307.2
$ws.Range("I126").Value = 10734.5
$ws.Range("J126").Value = 4909.1055
$ws.Range("K126").Value = 32203.5
$ws.Range("L126").Value = 14727.3165
$ws.Range("M126").Value = -29733.5
$ws.Range("N126").Value = -19667.3165
$ws.Range("H136").Value = 45456750
$ws.Range("I136").Value = 55558000
$ws.Range("J136").Value = 1150
$ws.Range("K136").Value = 166674000
$ws.Range("L136").Value = 3450
$ws.Range("M136").Value = -166671450
$ws.Range("N136").Value = -8550

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7070.857
$ws.Range("I62").Value = 6332.6665
$ws.Range("J62").Value = 7624.5
$ws.Range("K62").Value = 6332.6665
$ws.Range("L62").Value = 7624.5
$ws.Range("M62").Value = -5708.6665
$ws.Range("N62").Value = -8872.5
$ws.Range("H65").Value = 7070.857
$ws.Range("I65").Value = 6332.6665
$ws.Range("J65").Value = 7624.5
$ws.Range("K65").Value = 31663.3325
$ws.Range("L65").Value = 38122.5
$ws.Range("M65").Value = -28543.3325
$ws.Range("N65").Value = -44362.5
$ws.Range("H115").Value = 34663.332
$ws.Range("J115").Value = 34663.332
$ws.Range("L115").Value = 34663.332
$ws.Range("N115").Value = -37797.332
$ws.Range("H122").Value = 3982.6365
$ws.Range("I122").Value = 5139.6
$ws.Range("J122").Value = 3018.5
$ws.Range("K122").Value = 15418.8
$ws.Range("L122").Value = 9055.5
$ws.Range("M122").Value = -12968.8
$ws.Range("N122").Value = -13955.5
$ws.Range("H136").Value = 1065.0476
$ws.Range("I136").Value = 1018.3
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 3054.9
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -504.8999999999996
$ws.Range("N136").Value = -11100
